# Fix the trailing-space typo in the "abdominoperineal resection" search
# term on the English sheet so assertion failure messages read more
# cleanly (more descriptive AssertTrue() messages).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("SitewideSearchEn")
$ws.Range("A8").Value = "abdominoperineal resection"

# Restore focus to the English sheet (it is the first/active tab in the
# saved workbook).
$ws.Activate()
